# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-21 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 3
    4  = 6
    5  = 6
    6  = 5
    7  = 7
    8  = 3
    9  = 6
    10 = 7
    11 = 7
    12 = 4
    13 = 4
    14 = 4
    15 = 8
    16 = 3
    17 = 2
    18 = 5
    19 = 7
    20 = 4
    21 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
